$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constants (OLE COLORREF / BGR order, as used by Excel COM Font.Color)
$red   = 255      # 0x0000FF -> R=255,G=0,B=0   -> matches existing "done" style (red font)
$green = 32768    # 0x008000 -> R=0,G=128,B=0   -> matches existing "running"/"queueing" style (green font)

# Row 8: G8 "queueing (partial)" -> "done"; H8 newly filled with "running"
$ws.Range("G8").Value = "done"
$ws.Range("G8").Font.Color = $red
$ws.Range("H8").Value = "running"
$ws.Range("H8").Font.Color = $green

# Row 9: H9 "running" -> "done"
$ws.Range("H9").Value = "done"
$ws.Range("H9").Font.Color = $red

# Row 10: H10 "running" -> "done"
$ws.Range("H10").Value = "done"
$ws.Range("H10").Font.Color = $red

# Row 13: H13 "queueing" -> "done"
$ws.Range("H13").Value = "done"
$ws.Range("H13").Font.Color = $red

# Row 14: H14 "running" -> "done"
$ws.Range("H14").Value = "done"
$ws.Range("H14").Font.Color = $red

# Row 15: H15 "queueing" -> "running"
$ws.Range("H15").Value = "running"
$ws.Range("H15").Font.Color = $green

# Row 16: H16 "queueing" -> "running"
$ws.Range("H16").Value = "running"
$ws.Range("H16").Font.Color = $green

# Row 17: H17 "running" -> "done"
$ws.Range("H17").Value = "done"
$ws.Range("H17").Font.Color = $red

# Row 20: G20 "running" -> "done"; H20 newly filled with "done"
$ws.Range("G20").Value = "done"
$ws.Range("G20").Font.Color = $red
$ws.Range("H20").Value = "done"
$ws.Range("H20").Font.Color = $red

# Row 21: H21 "queueing" -> "running"
$ws.Range("H21").Value = "running"
$ws.Range("H21").Font.Color = $green

# Row 24: H24 "running" -> "done"
$ws.Range("H24").Value = "done"
$ws.Range("H24").Font.Color = $red
